# The final patch duplicates four entries from column A into a new column D
# (same row, same value) — rows 23, 33, 62 and 64 on the active sheet.
# Using Value2 (rather than Value) ensures the existing shared-string text is
# reused instead of being re-wrapped as a new distinct string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsToMirror = @(23, 33, 62, 64)

foreach ($r in $rowsToMirror) {
    $sourceCell = $ws.Cells.Item($r, 1)   # column A
    $targetCell = $ws.Cells.Item($r, 4)   # column D
    $targetCell.Value2 = $sourceCell.Value2
}
